# Applies the numeric corrections to the Leve profit tables across all sheets,
# matching the upstream "chore: update Sheets via scheduled runner" commit.
# Each block is guarded by a check on the row's Leve Item ID (column G) so a
# layout drift fails loudly instead of silently writing to the wrong row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$g = $ws.Range("G17").Value2
if ($g -ne 38956) { throw "ALC!G17: expected 38956, got $g" }
$ws.Range("H17").Value2 = 3000
$ws.Range("J17").Value2 = 3000
$ws.Range("L17").Value2 = 9000
$ws.Range("N17").Value2 = -9336

$g = $ws.Range("G47").Value2
if ($g -ne 2169) { throw "ALC!G47: expected 2169, got $g" }
$ws.Range("H47").Value2 = 7537
$ws.Range("I47").Value2 = 5000
$ws.Range("J47").Value2 = 10074
$ws.Range("K47").Value2 = 5000
$ws.Range("L47").Value2 = 10074
$ws.Range("M47").Value2 = -4028
$ws.Range("N47").Value2 = -12018

$ws = $wb.Worksheets.Item("ARM")
$g = $ws.Range("G97").Value2
if ($g -ne 19941) { throw "ARM!G97: expected 19941, got $g" }
$ws.Range("H97").Value2 = 66667370
$ws.Range("I97").Value2 = 66667370
$ws.Range("J97").Value2 = 0
$ws.Range("K97").Value2 = 66667370
$ws.Range("L97").Value2 = 0
$ws.Range("M97").Value2 = -66666874
$ws.Range("N97").ClearContents()

$g = $ws.Range("G132").Value2
if ($g -ne 43997) { throw "ARM!G132: expected 43997, got $g" }
$ws.Range("H132").Value2 = 1986.875
$ws.Range("I132").Value2 = 1986.875
$ws.Range("K132").Value2 = 5960.625
$ws.Range("M132").Value2 = -3430.625

$ws = $wb.Worksheets.Item("BSM")
$g = $ws.Range("G20").Value2
if ($g -ne 14149) { throw "BSM!G20: expected 14149, got $g" }
$ws.Range("H20").Value2 = 1070.25
$ws.Range("J20").Value2 = 1087.5
$ws.Range("L20").Value2 = 1087.5
$ws.Range("N20").Value2 = -1581.5

$g = $ws.Range("G64").Value2
if ($g -ne 14184) { throw "BSM!G64: expected 14184, got $g" }
$ws.Range("H64").Value2 = 1007
$ws.Range("J64").Value2 = 1007
$ws.Range("L64").Value2 = 1007
$ws.Range("N64").Value2 = -1457

$g = $ws.Range("G67").Value2
if ($g -ne 14184) { throw "BSM!G67: expected 14184, got $g" }
$ws.Range("H67").Value2 = 1007
$ws.Range("J67").Value2 = 1007
$ws.Range("L67").Value2 = 1007
$ws.Range("N67").Value2 = -2567

$g = $ws.Range("G70").Value2
if ($g -ne 15553) { throw "BSM!G70: expected 15553, got $g" }
$ws.Range("H70").Value2 = 245000
$ws.Range("J70").Value2 = 245000
$ws.Range("L70").Value2 = 245000
$ws.Range("N70").Value2 = -245586

$g = $ws.Range("G73").Value2
if ($g -ne 15553) { throw "BSM!G73: expected 15553, got $g" }
$ws.Range("H73").Value2 = 245000
$ws.Range("J73").Value2 = 245000
$ws.Range("L73").Value2 = 245000
$ws.Range("N73").Value2 = -247028

$ws = $wb.Worksheets.Item("CRP")
$g = $ws.Range("G26").Value2
if ($g -ne 2004) { throw "CRP!G26: expected 2004, got $g" }
$ws.Range("H26").Value2 = 650
$ws.Range("I26").Value2 = 650
$ws.Range("K26").Value2 = 650
$ws.Range("M26").Value2 = -363

$g = $ws.Range("G31").Value2
if ($g -ne 44023) { throw "CRP!G31: expected 44023, got $g" }
$ws.Range("H31").Value2 = 3589.6924
$ws.Range("I31").Value2 = 2308.6667
$ws.Range("J31").Value2 = 6472
$ws.Range("K31").Value2 = 2308.6667
$ws.Range("L31").Value2 = 6472
$ws.Range("M31").Value2 = -2013.6667
$ws.Range("N31").Value2 = -7062

$g = $ws.Range("G34").Value2
if ($g -ne 44023) { throw "CRP!G34: expected 44023, got $g" }
$ws.Range("H34").Value2 = 3589.6924
$ws.Range("I34").Value2 = 2308.6667
$ws.Range("J34").Value2 = 6472
$ws.Range("K34").Value2 = 2308.6667
$ws.Range("L34").Value2 = 6472
$ws.Range("M34").Value2 = -2106.6667
$ws.Range("N34").Value2 = -6876

$g = $ws.Range("G47").Value2
if ($g -ne 1920) { throw "CRP!G47: expected 1920, got $g" }
$ws.Range("H47").Value2 = 26124.875
$ws.Range("I47").Value2 = 40000
$ws.Range("J47").Value2 = 24142.715
$ws.Range("K47").Value2 = 40000
$ws.Range("L47").Value2 = 24142.715
$ws.Range("M47").Value2 = -39434
$ws.Range("N47").Value2 = -25274.715

$g = $ws.Range("G58").Value2
if ($g -ne 44021) { throw "CRP!G58: expected 44021, got $g" }
$ws.Range("H58").Value2 = 3012
$ws.Range("I58").Value2 = 3012
$ws.Range("K58").Value2 = 3012
$ws.Range("M58").Value2 = -2809

$g = $ws.Range("G59").Value2
if ($g -ne 1942) { throw "CRP!G59: expected 1942, got $g" }
$ws.Range("H59").Value2 = 60000
$ws.Range("I59").Value2 = 60000
$ws.Range("K59").Value2 = 60000
$ws.Range("M59").Value2 = -58855

$g = $ws.Range("G63").Value2
if ($g -ne 10604) { throw "CRP!G63: expected 10604, got $g" }
$ws.Range("H63").Value2 = 0
$ws.Range("J63").Value2 = 0
$ws.Range("L63").Value2 = 0
$ws.Range("N63").ClearContents()

$g = $ws.Range("G66").Value2
if ($g -ne 10604) { throw "CRP!G66: expected 10604, got $g" }
$ws.Range("H66").Value2 = 0
$ws.Range("J66").Value2 = 0
$ws.Range("L66").Value2 = 0
$ws.Range("N66").ClearContents()

$g = $ws.Range("G88").Value2
if ($g -ne 10608) { throw "CRP!G88: expected 10608, got $g" }
$ws.Range("H88").Value2 = 24009.916
$ws.Range("J88").Value2 = 24009.916
$ws.Range("L88").Value2 = 24009.916
$ws.Range("N88").Value2 = -24821.916

$g = $ws.Range("G91").Value2
if ($g -ne 10608) { throw "CRP!G91: expected 10608, got $g" }
$ws.Range("H91").Value2 = 24009.916
$ws.Range("J91").Value2 = 24009.916
$ws.Range("L91").Value2 = 24009.916
$ws.Range("N91").Value2 = -26817.916

$g = $ws.Range("G105").Value2
if ($g -ne 19928) { throw "CRP!G105: expected 19928, got $g" }
$ws.Range("H105").Value2 = 914
$ws.Range("I105").Value2 = 696.8
$ws.Range("K105").Value2 = 696.8
$ws.Range("M105").Value2 = 1050.2

$g = $ws.Range("G132").Value2
if ($g -ne 44019) { throw "CRP!G132: expected 44019, got $g" }
$ws.Range("H132").Value2 = 2203.4443
$ws.Range("I132").Value2 = 1852.125
$ws.Range("J132").Value2 = 5014
$ws.Range("K132").Value2 = 5556.375
$ws.Range("L132").Value2 = 15042
$ws.Range("M132").Value2 = -3026.375
$ws.Range("N132").Value2 = -20102

$g = $ws.Range("G136").Value2
if ($g -ne 44021) { throw "CRP!G136: expected 44021, got $g" }
$ws.Range("H136").Value2 = 3012
$ws.Range("I136").Value2 = 3012
$ws.Range("K136").Value2 = 9036
$ws.Range("M136").Value2 = -6486

$ws = $wb.Worksheets.Item("CUL")
$g = $ws.Range("G106").Value2
if ($g -ne 19819) { throw "CUL!G106: expected 19819, got $g" }
$ws.Range("H106").Value2 = 0
$ws.Range("J106").Value2 = 0
$ws.Range("L106").Value2 = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$g = $ws.Range("G49").Value2
if ($g -ne 4232) { throw "GSM!G49: expected 4232, got $g" }
$ws.Range("H49").Value2 = 25927.6
$ws.Range("J49").Value2 = 29909.5
$ws.Range("L49").Value2 = 29909.5
$ws.Range("N49").Value2 = -30277.5

$g = $ws.Range("G122").Value2
if ($g -ne 36182) { throw "GSM!G122: expected 36182, got $g" }
$ws.Range("H122").Value2 = 5036.6
$ws.Range("I122").Value2 = 3481.6428
$ws.Range("K122").Value2 = 10444.9284
$ws.Range("M122").Value2 = -7994.928400000001

$ws = $wb.Worksheets.Item("LTW")
$g = $ws.Range("G68").Value2
if ($g -ne 12563) { throw "LTW!G68: expected 12563, got $g" }
$ws.Range("H68").Value2 = 6232.5
$ws.Range("J68").Value2 = 4133.3335
$ws.Range("L68").Value2 = 4133.3335
$ws.Range("N68").Value2 = -5631.3335

$g = $ws.Range("G71").Value2
if ($g -ne 12563) { throw "LTW!G71: expected 12563, got $g" }
$ws.Range("H71").Value2 = 6232.5
$ws.Range("J71").Value2 = 4133.3335
$ws.Range("L71").Value2 = 20666.6675
$ws.Range("N71").Value2 = -28154.6675

$g = $ws.Range("G82").Value2
if ($g -ne 12565) { throw "LTW!G82: expected 12565, got $g" }
$ws.Range("H82").Value2 = 3650
$ws.Range("I82").Value2 = 1666.6666
$ws.Range("J82").Value2 = 4500
$ws.Range("K82").Value2 = 1666.6666
$ws.Range("L82").Value2 = 4500
$ws.Range("M82").Value2 = -1305.6666
$ws.Range("N82").Value2 = -5222

$g = $ws.Range("G85").Value2
if ($g -ne 12565) { throw "LTW!G85: expected 12565, got $g" }
$ws.Range("H85").Value2 = 3650
$ws.Range("I85").Value2 = 1666.6666
$ws.Range("J85").Value2 = 4500
$ws.Range("K85").Value2 = 1666.6666
$ws.Range("L85").Value2 = 4500
$ws.Range("M85").Value2 = -418.6666
$ws.Range("N85").Value2 = -6996

$g = $ws.Range("G136").Value2
if ($g -ne 44060) { throw "LTW!G136: expected 44060, got $g" }
$ws.Range("H136").Value2 = 5000
$ws.Range("I136").Value2 = 5000
$ws.Range("K136").Value2 = 15000
$ws.Range("M136").Value2 = -12450

$ws = $wb.Worksheets.Item("WVR")
$g = $ws.Range("G48").Value2
if ($g -ne 3140) { throw "WVR!G48: expected 3140, got $g" }
$ws.Range("H48").Value2 = 0
$ws.Range("J48").Value2 = 0
$ws.Range("L48").Value2 = 0
$ws.Range("N48").ClearContents()

$g = $ws.Range("G62").Value2
if ($g -ne 12589) { throw "WVR!G62: expected 12589, got $g" }
$ws.Range("H62").Value2 = 3699.75
$ws.Range("I62").Value2 = 3266.3333
$ws.Range("J62").Value2 = 5000
$ws.Range("K62").Value2 = 3266.3333
$ws.Range("L62").Value2 = 5000
$ws.Range("M62").Value2 = -2642.3333
$ws.Range("N62").Value2 = -6248

$g = $ws.Range("G65").Value2
if ($g -ne 12589) { throw "WVR!G65: expected 12589, got $g" }
$ws.Range("H65").Value2 = 3699.75
$ws.Range("I65").Value2 = 3266.3333
$ws.Range("J65").Value2 = 5000
$ws.Range("K65").Value2 = 16331.6665
$ws.Range("L65").Value2 = 25000
$ws.Range("M65").Value2 = -13211.6665
$ws.Range("N65").Value2 = -31240

$g = $ws.Range("G107").Value2
if ($g -ne 27746) { throw "WVR!G107: expected 27746, got $g" }
$ws.Range("H107").Value2 = 288.75
$ws.Range("I107").Value2 = 233.33333
$ws.Range("J107").Value2 = 455
$ws.Range("K107").Value2 = 699.99999
$ws.Range("L107").Value2 = 1365
$ws.Range("M107").Value2 = 1220.00001
$ws.Range("N107").Value2 = -5205

$g = $ws.Range("G132").Value2
if ($g -ne 44029) { throw "WVR!G132: expected 44029, got $g" }
$ws.Range("H132").Value2 = 2250
$ws.Range("I132").Value2 = 2250
$ws.Range("K132").Value2 = 6750
$ws.Range("M132").Value2 = -4220
